# Nacelle Mass Properties: fix rounding format error
# - "Mass" column (D) was recorded in tonnes with a 3-decimal custom format;
#   convert the stored values to kilograms (x1000) and switch the column to
#   a whole-number "0" format so kg values don't show spurious decimals.
# - Update the "Units" label for that column from "t" to "kg".
# - Flip the sign of the X_TT (column B) coordinates for all component rows
#   except the ones that are already zero.
# - Restore the active selection/cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Units label: t -> kg ----
$ws.Range("D2").Value = "kg"

# ---- Mass column (D): tonnes -> kilograms, integer display format ----
$massUpdates = @{
    "D3"  = 100000
    "D4"  = 13362.348
    "D5"  = 226700
    "D6"  = 145250
    "D7"  = 19503.66
    "D8"  = 190000
    "D9"  = 39433.565000000002
    "D10" = 3627.46
    "D11" = 50000
    "D12" = 4699
    "D13" = 4699
    "D14" = 797275.03899999999
}

foreach ($addr in $massUpdates.Keys) {
    $ws.Range($addr).Value = $massUpdates[$addr]
    $ws.Range($addr).NumberFormat = "0"
}

# ---- X_TT column (B): flip sign (rows whose value is 0 stay untouched) ----
$signFlips = @{
    "B4"  = -6.0750000000000002
    "B5"  = -5.7350000000000003
    "B6"  = -6.3390000000000004
    "B7"  = -6.1390000000000002
    "B8"  = -10.685
    "B9"  = -0.80900000000000005
    "B10" = -4.6710000000000003
    "B12" = -5.3559999999999999
    "B13" = -6.9480000000000004
    "B14" = -5.718
}

foreach ($addr in $signFlips.Keys) {
    $ws.Range($addr).Value = $signFlips[$addr]
}

# ---- Restore selection ----
$ws.Range("E16").Select()
